$wb = $excel.ActiveWorkbook

$wsFile = $wb.Worksheets.Item("FILESTREAM")

$wsFile.Range("A2").Value = "EXAMPLE"
$wsFile.Range("E2").Value = "example@example.com"

$wsFile.Select()
$wsFile.Range("A2").Select()
